$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D sometimes holds numeric-looking text (e.g. "223.58", "1.00", "0.630").
# Plain `.Value` assignment would let Excel auto-coerce those into real numbers,
# which would change their cell type/formatting vs. the original inline-string
# cells. Force them to stay text by toggling NumberFormat to Text while assigning
# the value, then clear the format again so no stray style survives the edit.
function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.ClearFormats()
}

# Row 2
Set-TextValue $ws.Range("D2") '34.626.63'
$ws.Range("E2").Value = '  +2.60%  '

# Row 3
Set-TextValue $ws.Range("D3") '1.788.75'
$ws.Range("E3").Value = '  +0.66%  '

# Row 4
$ws.Range("E4").Value = '  +0.03%  '

# Row 5
Set-TextValue $ws.Range("D5") '223.58'
$ws.Range("E5").Value = '  -0.66%  '

# Row 6
Set-TextValue $ws.Range("D6") '0.556'
$ws.Range("E6").Value = '  -0.29%  '

# Row 7
$ws.Range("E7").Value = '  +0.05%  '

# Row 8
Set-TextValue $ws.Range("D8") '33.04'
$ws.Range("E8").Value = '  +7.69%  '

# Row 9
$ws.Range("E9").Value = '  +0.99%  '

# Row 10
Set-TextValue $ws.Range("D10") '0.0678'
$ws.Range("E10").Value = '  +2.59%  '

# Row 11
$ws.Range("E11").Value = '  +1.40%  '

# Row 12
Set-TextValue $ws.Range("D12") '2.045.22'
$ws.Range("E12").Value = '  +0.75%  '

# Row 13
$ws.Range("B13").Value = 'WrappedEther'
$ws.Range("C13").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
Set-TextValue $ws.Range("D13") '1.826.61'
$ws.Range("E13").Value = '  +2.83%  '

# Row 14
$ws.Range("B14").Value = 'Chainlink'
$ws.Range("C14").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
Set-TextValue $ws.Range("D14") '11.15'
$ws.Range("E14").Value = '  +11.62%  '

# Row 15
Set-TextValue $ws.Range("D15") '34.612.59'
$ws.Range("E15").Value = '  +2.56%  '

# Row 16
Set-TextValue $ws.Range("D16") '0.631'
$ws.Range("E16").Value = '  +0.44%  '

# Row 17
Set-TextValue $ws.Range("D17") '4.30'
$ws.Range("E17").Value = '  +2.84%  '

# Row 18
Set-TextValue $ws.Range("D18") '68.53'
$ws.Range("E18").Value = '  +0.02%  '

# Row 19
Set-TextValue $ws.Range("D19") '253.34'
$ws.Range("E19").Value = '  +0.66%  '

# Row 20
Set-TextValue $ws.Range("D20") '0.0₃0772'
$ws.Range("E20").Value = '  +4.72%  '

# Row 21
$ws.Range("E21").Value = '  +0.12%  '

# Row 22
$ws.Range("E22").Value = '  +1.56%  '

# Row 23
$ws.Range("E23").Value = '  +1.14%  '

# Row 24
Set-TextValue $ws.Range("D24") '2.15'
$ws.Range("E24").Value = '  +0.29%  '

# Row 25
Set-TextValue $ws.Range("D25") '158.51'
$ws.Range("E25").Value = '  -0.29%  '

# Row 26
Set-TextValue $ws.Range("D26") '16.32'
$ws.Range("E26").Value = '  -1.09%  '

# Row 27
Set-TextValue $ws.Range("D27") '7.09'
$ws.Range("E27").Value = '  +2.32%  '

# Row 28
$ws.Range("E28").Value = '  -0.11%  '

# Row 29
Set-TextValue $ws.Range("D29") '1.00'
$ws.Range("E29").Value = '  -0.08%  '

# Row 30
$ws.Range("B30").Value = 'Filecoin'
$ws.Range("C30").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
Set-TextValue $ws.Range("D30") '3.75'
$ws.Range("E30").Value = '  -1.50%  '

# Row 31
$ws.Range("B31").Value = 'Hedera'
$ws.Range("C31").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
Set-TextValue $ws.Range("D31") '0.0514'
$ws.Range("E31").Value = '  +0.18%  '

# Row 32
$ws.Range("E32").Value = '  +0.26%  '

# Row 33
$ws.Range("E33").Value = '  +0.52%  '

# Row 34
$ws.Range("E34").Value = '  +1.58%  '

# Row 35
Set-TextValue $ws.Range("D35") '1.442.65'
$ws.Range("E35").Value = '  -2.75%  '

# Row 36
$ws.Range("E36").Value = '  -0.84%  '

# Row 37
$ws.Range("B37").Value = 'VeChain'
$ws.Range("C37").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
Set-TextValue $ws.Range("D37") '0.0189'
$ws.Range("E37").Value = '  +2.13%  '

# Row 38
$ws.Range("B38").Value = 'ImmutableX'
$ws.Range("C38").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
Set-TextValue $ws.Range("D38") '0.630'
$ws.Range("E38").Value = '  -0.70%  '

# Row 39
Set-TextValue $ws.Range("D39") '83.15'
$ws.Range("E39").Value = '  -0.19%  '

# Row 40
$ws.Range("E40").Value = '  +3.95%  '

# Row 41
$ws.Range("E41").Value = '  +0.78%  '

# Row 42
Set-TextValue $ws.Range("D42") '0.902'
$ws.Range("E42").Value = '  +1.94%  '

# Row 43
$ws.Range("E43").Value = '  -0.90%  '

# Row 44
Set-TextValue $ws.Range("D44") '0.0504'
$ws.Range("E44").Value = '  -1.71%  '

# Row 45
$ws.Range("E45").Value = '  +2.46%  '

# Row 46
$ws.Range("E46").Value = '  -2.18%  '

# Row 47
Set-TextValue $ws.Range("D47") '1.942.98'
$ws.Range("E47").Value = '  +0.65%  '

# Row 48
Set-TextValue $ws.Range("D48") '104.69'
$ws.Range("E48").Value = '  +7.18%  '

# Row 49
$ws.Range("B49").Value = 'InjectiveProtocol'
$ws.Range("C49").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
Set-TextValue $ws.Range("D49") '11.99'
$ws.Range("E49").Value = '  +2.83%  '

# Row 50
$ws.Range("B50").Value = 'PaxDollar'
$ws.Range("C50").Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
Set-TextValue $ws.Range("D50") '1.00'
$ws.Range("E50").Value = '  +0.06%  '

# Row 51
Set-TextValue $ws.Range("D51") '49.29'
$ws.Range("E51").Value = '  -2.87%  '
